$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column R holds the 2021 figures, matching the layout/styling already
# used by the other year columns (D:Q). Copy formats from the nearest
# existing cell that carries the right look for each row:
#  - row 4 (year header)            -> same style as Q4
#  - row 5 (share of renewables, %) -> same style as D5 (no extra decimal fmt)
#  - row 6 (hydropower production)  -> same style as Q6
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D5").Copy()
$ws.Range("R5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("Q6").Copy()
$ws.Range("R6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Add new column R data for year 2021
$ws.Range("R4").Value = 2021
$ws.Range("R5").Value = 31.8
$ws.Range("R6").Value = 12957.1

# Update selection to match target state
$ws.Range("R4:R6").Select()
